$wb = $excel.ActiveWorkbook

$oldText = "Subject completes response to perturbation having steered the vehicle back to the center of the lane. Normally this would be tagged with temporal scope, but avoiding definitions here."
$newText = "Subject completes response to perturbation having steered the vehicle back to the center of the lane. Normally this would be tagged with temporal scope but avoiding definitions here."

# Correct the typo (stray comma) in the "LKT 8HED3" sheet
$ws1 = $wb.Worksheets.Item("LKT 8HED3")
if ($ws1.Range("D5").Value2 -eq $oldText) {
    $ws1.Range("D5").Value2 = $newText
}

# Correct the same typo in the "LKT 8HED3A" sheet
$ws5 = $wb.Worksheets.Item("LKT 8HED3A")
if ($ws5.Range("D5").Value2 -eq $oldText) {
    $ws5.Range("D5").Value2 = $newText
}

# Reflect the final selection / active sheet state left behind by the edit
$ws1.Range("D5").Select() | Out-Null
$ws5.Range("D5").Select() | Out-Null

$ws4 = $wb.Worksheets.Item("DAS Events")
$ws4.Activate() | Out-Null
